$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Hammoud, Helal Mustafa) was moved from the Planning group/team (LP*)
# to the Corporate IT group / Services & Call Center team (LT*), and is now
# flagged as a Team Leader.

# ORG_CODE: LP21 -> LT51
$ws.Range("C5").Value = "LT51"

# GROUP_Code: LP01 -> LT01 ; remove the old top border highlight on this cell
$ws.Range("K5").Value = "LT01"
$ws.Range("K5").Borders.Item(8).LineStyle = -4142

# GROUP_NAME: LP01-Planning Group -> LT01-CORPORATE INFORMATION TECHNOLOGY GROUP.
$ws.Range("L5").Value = "LT01-CORPORATE INFORMATION TECHNOLOGY GROUP."
$ws.Range("L5").Style = "Normal"
$ws.Range("L5").Font.Color = 0
$ws.Range("L5").Borders.Item(8).Color = 13999631
$ws.Range("L5").Borders.Item(8).LineStyle = 1
$ws.Range("L5").Borders.Item(9).Color = 13999631
$ws.Range("L5").Borders.Item(9).LineStyle = 1

# TEAM_Code: LP11 -> LT51
$ws.Range("O5").Value = "LT51"

# TEAM_NAME: LP11-Strategic Planning Team -> LT51 - Services & Call Center Team
# remove the old top border highlight on this cell
$ws.Range("P5").Value = "LT51 - Services & Call Center Team"
$ws.Range("P5").Borders.Item(8).LineStyle = -4142

# Is_TeamLeader: 0 -> 1
$ws.Range("S5").Value = 1

# Update the active selection left over from the edit session
$ws.Range("A5").Select()
